# Se ajuste foco aplicación InfoNova a traves del selector
# Adds a new "Rango Final" column (F) and updates the first result row
# to reflect the new successful transaction.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column F
$ws.Range("F1").Value = "Rango Final"

# Update the first data row (row 2) to the new successful transaction values
$ws.Range("B2").Value = "00837202410161513058-18102024"
$ws.Range("C2").Value = "OK"
$ws.Range("D2").Value = "Successful Transaction: 1 - Reference: 00837202410161513058-18102024"
$ws.Range("E2").Value = "TransactionTime: 0h 1m 36s"

# The new "Rango Final" value is numeric-looking but stored as text
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "948690199"
$ws.Range("F2").Style = "Normal"
